# CATATAN REVISI.xlsx - "update CATATAN REVISI.xlsx dan data baru dari pa devi"
#
# Fills in rows 106-112 of Sheet1 with new revision notes / statuses that
# were supplied by Pa Devi, renumbers column A (NO) for the newly
# filled-in rows, and moves the viewport/selection to reflect where the
# user ended up working (around row 93-113).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 106: existing note just got marked "ok" ---
$ws.Range("C106").Value = "ok"

# --- Row 107: new note, no status yet ---
$ws.Range("A107").Value = 135
$ws.Range("B107").Value = "bisa cetak surat usulan, file rincian kegiatan (semua)"

# --- Row 108: new note, marked "ok" ---
$ws.Range("A108").Value = 136
$ws.Range("B108").Value = "di umkm: ditampilkan nama, hapus kelurahan, kecamatan, kota"
$ws.Range("C108").Value = "ok"

# --- Row 109: new note, no status yet ---
$ws.Range("A109").Value = 137
$ws.Range("B109").Value = "gallery di link sama hasil2 pembangunan"

# --- Row 110: new note, no status yet ---
$ws.Range("A110").Value = 138
$ws.Range("B110").Value = "perbaikan bug foto dan dokumen"

# --- Row 111: new note, no status yet ---
$ws.Range("A111").Value = 139
$ws.Range("B111").Value = "database di hapus, diganti dengan data2 yang baru"

# --- Row 112: new note, marked "ok" (no NO value for this one) ---
$ws.Range("B112").Value = "kata-kata di bawah blog di hapus"
$ws.Range("C112").Value = "ok"

# --- Update the viewport / current selection to where the edits left off ---
$win = $excel.ActiveWindow
$win.ScrollRow = 93
$win.ScrollColumn = 1
$ws.Range("B113").Select()
